$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$experience = @(
    "Experience",
    "3 - 8 Years",
    "10 - 14 Years",
    "4 - 8 Years",
    "8 - 12 Years",
    "4 - 8 Years",
    "3 - 10 Years",
    "3 - 5 Years",
    "3 - 8 Years",
    "3 - 6 Years",
    "10 - 14 Years",
    "4 - 8 Years",
    "3 - 10 Years",
    "2 - 5 Years",
    "3 - 8 Years",
    "3 - 8 Years",
    "3 - 8 Years",
    "4 - 8 Years",
    "8 - 14 Years",
    "4 - 8 Years",
    "3 - 8 Years",
    "3 - 6 Years",
    "3 - 6 Years",
    "3 - 8 Years",
    "3 - 8 Years"
)

$status = @(
    "Status",
    "Open Position",
    "Open Position",
    "Open Position",
    "Open Position",
    "Open Position",
    "Open Position",
    "Open Position",
    "Open Position",
    "Open Position",
    "Open Position",
    "Open Position",
    "Open Position",
    "Open Position",
    "Open Position",
    "Open Position",
    "Open Position",
    "Open Position",
    "Open Position",
    "Open Position",
    "Open Position",
    "Open Position",
    "Open Position",
    "Open Position",
    "Open Position"
)

for ($i = 0; $i -lt 25; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 3).Value = $experience[$i]
    $ws.Cells.Item($row, 4).Value = $status[$i]
}
